$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- I column: AVERAGE formulas, J column: GECTI/KALDI result text ---

$ws.Range("I4").Formula = "=AVERAGE(D3:G4)"
$ws.Range("J4").Value = "GEÇTİ"

$ws.Range("I5").Formula = "=AVERAGE(D4:G5)"
$ws.Range("J5").Value = "GEÇTİ"

$ws.Range("I6").Formula = "=AVERAGE(D5:G6)"
$ws.Range("J6").Value = "GEÇTİ"

$ws.Range("I7").Formula = "=AVERAGE(D6:G7)"
$ws.Range("J7").Value = "GEÇTİ"

$ws.Range("I8").Formula = "=AVERAGE(D7:G8)"
$ws.Range("J8").Value = "KALDI"

$ws.Range("I9").Formula = "=AVERAGE(D8:G9)"
$ws.Range("J9").Value = "GEÇTİ"

$ws.Range("I10").Formula = "=AVERAGE(D9:G10)"
$ws.Range("J10").Value = "GEÇTİ"

$ws.Range("I11").Formula = "=AVERAGE(D10:G11)"
$ws.Range("J11").Value = "GEÇTİ"

# --- Row 12: class total scores ---
$ws.Range("D12").Value = 510
$ws.Range("E12").Value = 522
$ws.Range("F12").Value = 535
$ws.Range("G12").Value = 527
$ws.Range("H12").Formula = "=SUM(D12:G12)"

# --- Row 13: class averages ---
$ws.Range("D13").Formula = "=AVERAGE(D4:D11)"
$ws.Range("E13").Formula = "=AVERAGE(E4:E11)"
$ws.Range("F13").Formula = "=AVERAGE(F4:F11)"
$ws.Range("G13").Formula = "=AVERAGE(G4:G11)"

# --- Student info block ---
$ws.Range("F17").Value = 20215070055
$ws.Range("F18").Value = "Muhammed Ali Harmancı"
$ws.Range("F19").Value = "Yönetim Bilişim Sistemleri"

# --- Formatting tweaks ---
$ws.Rows.Item(1).RowHeight = 23.3
$ws.Rows.Item(11).RowHeight = 15.8
$ws.Rows.Item(13).RowHeight = 21.1
$ws.Rows.Item(14).RowHeight = 18.35

$ws.Columns.Item(8).ColumnWidth = 0.125
$ws.Columns.Item(9).ColumnWidth = 20.125
$ws.Columns.Item(10).ColumnWidth = 10.5

$ws.Range("F18:H18").Select()

$wb.Save()
